$d = $word.ActiveDocument

# --- Split the run containing "{m" into two runs: "{" and "m" ---
$found1 = $d.Content
$found1.Find.Execute("{m", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$brace = $d.Range($found1.Start, $found1.Start + 1)
# Toggling a character-formatting property on just the leading character
# forces Word to split the run at that boundary; reverting the value
# right back keeps both resulting runs identically formatted.
$brace.Font.Bold = 1
$brace.Font.Bold = 0

# --- Split the run containing ")}" into two runs: ")" and "}" ---
$found2 = $d.Content
$found2.Find.Execute(")}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$closeParen = $d.Range($found2.Start, $found2.Start + 1)
$closeParen.Font.Bold = 1
$closeParen.Font.Bold = 0

# The trailing "}" becomes its own, unformatted run (no rPr copied from
# its neighbour): delete it and retype it as fresh text.
$found3 = $d.Content
$found3.Find.Execute("}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$closingBrace = $d.Range($found3.Start, $found3.Start + 1)
$closingBrace.Delete()
$insertionPoint = $d.Range($found3.Start, $found3.Start)
$insertionPoint.InsertAfter("}")
